$wb = $excel.ActiveWorkbook

$oldNames = @(
  "summ23561050","summ23772948","summ23992773","summ24307249","summ24538281",
  "summ24769961","summ24983987","summ25196026","summ25411513","summ25636391",
  "summ25878982","summ26098536","summ26304563","summ26525598","summ26756650",
  "summ26965717","summ27180015","summ27389038","summ27601066","summ27829468",
  "summ28063503","summ28277579","summ28484610","summ28692358","summ28901388",
  "summ29120311","summ29336784","summ29552338","summ29766892","summ29985917",
  "summ30202978","summ30418006","summ30632029","summ30850073","summ31091698",
  "summ31339736","summ31584818","summ31840210","summ32071593","summ32287633",
  "summ32501193","summ32716222","summ32933436","summ33155554","summ33368085",
  "summ33586870","summ33805567","summ34024974","summ34243006","summ34456034"
)

$newNames = @(
  "summ59533399","summ59733421","summ59981475","summ00374905","summ00639527",
  "summ00898632","summ01156689","summ01410682","summ01694899","summ01959828",
  "summ02238718","summ02499740","summ02761799","summ03021213","summ03328539",
  "summ03605118","summ03867184","summ04169025","summ04540287","summ04851388",
  "summ05142581","summ05451431","summ05724571","summ05975544","summ06216587",
  "summ06475126","summ06758452","summ07028507","summ07284034","summ07544600",
  "summ07796845","summ08105238","summ08364813","summ08609844","summ08866718",
  "summ09122425","summ09371583","summ09714138","summ09966435","summ10255067",
  "summ10523722","summ10776004","summ11025314","summ11294242","summ11540940",
  "summ11820431","summ12095818","summ12367703","summ12617264","summ12895321"
)

for ($i = 0; $i -lt $oldNames.Length; $i++) {
    $ws = $wb.Worksheets.Item($oldNames[$i])
    $ws.Name = $newNames[$i]
}
